$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row for the new employee (John Smith) while sheet still uses the old A-D layout
$ws.Rows("2").Insert()
$ws.Range("A2").Value = 'John'
$ws.Range("B2").Value = 'Smith'
$ws.Range("C2").Value = 'testpassword123'
$ws.Range("D2").Value = 'driver'

# 2) Shift employee data from columns A-D into B-E to make room for a new Employee ID column in A
#    (cell values are written directly so column-width metadata for columns B/C is left untouched,
#     matching how the source workbook records this edit)
$ws.Range("E1").Value = 'admin'
$ws.Range("D1").Value = 'password123'
$ws.Range("C1").Value = 'Doe'
$ws.Range("B1").Value = 'John'
$ws.Range("E2").Value = 'driver'
$ws.Range("D2").Value = 'testpassword123'
$ws.Range("C2").Value = 'Smith'
$ws.Range("B2").Value = 'John'
$ws.Range("E3").Value = 'driver'
$ws.Range("D3").Value = 'passpopek2'
$ws.Range("C3").Value = 'Popek'
$ws.Range("B3").Value = 'Sean'
$ws.Range("E4").Value = 'driver'
$ws.Range("D4").Value = 'royspass56'
$ws.Range("C4").Value = 'Lopiccolo'
$ws.Range("B4").Value = 'Roy '
$ws.Range("E5").Value = 'driver'
$ws.Range("D5").Value = 'gmirando93'
$ws.Range("C5").Value = 'Mirando'
$ws.Range("B5").Value = 'Gil'
$ws.Range("E6").Value = 'driver'
$ws.Range("D6").Value = '82karyPam'
$ws.Range("C6").Value = 'Kary'
$ws.Range("B6").Value = 'Pam'
$ws.Range("E7").Value = 'driver'
$ws.Range("D7").Value = '1978gconstant'
$ws.Range("C7").Value = 'Constant'
$ws.Range("B7").Value = 'Graig'
$ws.Range("E8").Value = 'driver'
$ws.Range("D8").Value = 'nohillberry1'
$ws.Range("C8").Value = 'Berryhill'
$ws.Range("B8").Value = 'Noble'
$ws.Range("E9").Value = 'driver'
$ws.Range("D9").Value = 'elephant381'
$ws.Range("C9").Value = 'Usher'
$ws.Range("B9").Value = 'Leigh'
$ws.Range("E10").Value = 'driver'
$ws.Range("D10").Value = 'treesandbees23'
$ws.Range("C10").Value = 'Snowden'
$ws.Range("B10").Value = 'Normand'
$ws.Range("E11").Value = 'admin'
$ws.Range("D11").Value = 'ferrari950'
$ws.Range("C11").Value = 'Molitor'
$ws.Range("B11").Value = 'Normand'
$ws.Range("E12").Value = 'driver'
$ws.Range("D12").Value = 'coachdriver837'
$ws.Range("C12").Value = 'Velarde'
$ws.Range("B12").Value = 'Ethelyn'
$ws.Range("E13").Value = 'driver'
$ws.Range("D13").Value = 'piestiesflies1'
$ws.Range("C13").Value = 'Matte'
$ws.Range("B13").Value = 'Judson'
$ws.Range("E14").Value = 'driver'
$ws.Range("D14").Value = 'wordpass321'
$ws.Range("C14").Value = 'Ahmad'
$ws.Range("B14").Value = 'Gus'
$ws.Range("E15").Value = 'driver'
$ws.Range("D15").Value = 'keyboardhello'
$ws.Range("C15").Value = 'Bodin'
$ws.Range("B15").Value = 'Mitchel'
$ws.Range("E16").Value = 'admin'
$ws.Range("D16").Value = 'magiccar11'
$ws.Range("C16").Value = 'Lagarde'
$ws.Range("B16").Value = 'Dane'
$ws.Range("E17").Value = 'admin'
$ws.Range("D17").Value = 'smallbigfig25'
$ws.Range("C17").Value = 'Marques'
$ws.Range("B17").Value = 'Jone'
$ws.Range("E18").Value = 'admin'
$ws.Range("D18").Value = 'big100money'
$ws.Range("C18").Value = 'Danna'
$ws.Range("B18").Value = 'Joi'
$ws.Range("E19").Value = 'driver'
$ws.Range("D19").Value = 'buildings22trees'
$ws.Range("C19").Value = 'Hammet'
$ws.Range("B19").Value = 'Arlene'
$ws.Range("E20").Value = 'driver'
$ws.Range("D20").Value = 'windows10frame'
$ws.Range("C20").Value = 'Siewert'
$ws.Range("B20").Value = 'Leah'
$ws.Range("E21").Value = 'driver'
$ws.Range("D21").Value = '1iconsadd89'
$ws.Range("C21").Value = 'Williams'
$ws.Range("B21").Value = 'Michele'

# 3) Fill in the new Employee ID column (A), top to bottom
$ws.Range("A1").Value = 'A1019'
$ws.Range("A2").Value = 'D1006'
$ws.Range("A3").Value = 'D1020'
$ws.Range("A4").Value = 'D1021'
$ws.Range("A5").Value = 'D1022'
$ws.Range("A6").Value = 'D1023'
$ws.Range("A7").Value = 'D1024'
$ws.Range("A8").Value = 'D1025'
$ws.Range("A9").Value = 'D1026'
$ws.Range("A10").Value = 'D1027'
$ws.Range("A11").Value = 'A1028'
$ws.Range("A12").Value = 'D1029'
$ws.Range("A13").Value = 'D1030'
$ws.Range("A14").Value = 'D1031'
$ws.Range("A15").Value = 'D1032'
$ws.Range("A16").Value = 'A1033'
$ws.Range("A17").Value = 'A1034'
$ws.Range("A18").Value = 'A1035'
$ws.Range("A19").Value = 'D1036'
$ws.Range("A20").Value = 'D1037'
$ws.Range("A21").Value = 'D1038'

# 4) New Password column width + restore selection
$ws.Columns("D").ColumnWidth = 18.6
$ws.Range("A22").Select() | Out-Null
